$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header labels in row 1: "testResultActions" -> "testResultDetails" ---
$ws.Range("A1").Value = "button_testResultDetails_class"
$ws.Range("B1").Value = "button_testResultDetails_class_1"
$ws.Range("C1").Value = "button_testResultDetails_class_2"
$ws.Range("D1").Value = "button_testResultDetails_class_3"
$ws.Range("E1").Value = "button_testResultDetails_class_4"
$ws.Range("F1").Value = "button_testResultDetails_class_5"
$ws.Range("G1").Value = "button_testResultDetails_internalRoleButtonName"
$ws.Range("H1").Value = "button_testResultDetails_internalRoleButtonName_1"
$ws.Range("I1").Value = "button_testResultDetails_internalRoleButtonName_2"
$ws.Range("J1").Value = "button_testResultDetails_internalRoleButtonName_3"

# --- Update data path values in row 2 ---
$ws.Range("K2").Value = "Data Files/AI-Generated/Common/viewTestResultsAndMarkAsPassed-test-data"
$ws.Range("L2").Value = "Data Files/AI-Generated/Common/viewTestResultsAndMarkAsPassed-test-data"

# --- Widen columns K (11) and L (12) from 65 to 73 ---
# Note: Excel's ColumnWidth property (character units) is offset from the
# stored OOXML "width" attribute by ~5/6 (the default column padding), so we
# subtract that offset here to land on an exact stored width of 73.
$targetWidth = 73 - (5/6)
$ws.Columns.Item(11).ColumnWidth = $targetWidth
$ws.Columns.Item(12).ColumnWidth = $targetWidth
